$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 543
$ws1.Range("F14").Value = 795
$ws1.Range("F15").Value = 6458
$ws1.Range("F19").Value = 4295
$ws1.Range("F23").Value = 4059
$ws1.Range("F25").Value = 200
$ws1.Range("F35").Value = 7295
$ws1.Range("F44").Value = 800
$ws1.Range("F46").Value = 3500

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 15

# Sheet "全部类型" (all types, aggregated view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 543
$ws4.Range("F15").Value = 795
$ws4.Range("F16").Value = 6458
$ws4.Range("F20").Value = 4295
$ws4.Range("F24").Value = 4059
$ws4.Range("F26").Value = 200
$ws4.Range("F33").Value = 7295
$ws4.Range("F43").Value = 800
$ws4.Range("F45").Value = 3500
